$d = $word.ActiveDocument

# The feedback table's "Answers" column (2nd column) contains one-word
# answers that need to be prefixed with "Design: ". Several of those
# words (e.g. "jäi") also occur as whole words elsewhere in the document
# (e.g. inside the matching Question cell), so a plain document-wide
# Find/Replace would touch the wrong occurrence. Instead, gather the
# exact answer cells first (row order == document order), then replace
# each one by searching forward from that cell's own start position,
# keeping a running offset "shift" to account for the text already
# inserted by earlier replacements.

$tbl = $d.Tables.Item(1)

$targets = @()
foreach ($cell in $tbl.Range.Cells) {
    if ($cell.ColumnIndex -eq 2 -and $cell.RowIndex -gt 1) {
        $wordText = $cell.Range.Text.Trim([char]13, [char]7)
        $targets += , @{ LowerBound = $cell.Range.Start; Word = $wordText }
    }
}

$shift = 0
foreach ($t in $targets) {
    $lowerBound = $t.LowerBound + $shift
    $searchRange = $d.Range($lowerBound, $d.Content.End)
    $found = $searchRange.Find.Execute($t.Word, $true, $true)
    if ($found) {
        $newText = "Design: " + $t.Word
        $oldLen = $searchRange.End - $searchRange.Start
        $searchRange.Text = $newText
        $shift += ($newText.Length - $oldLen)
    }
}
